$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05929980902570103
$ws.Range("H2").Value = -8.419554475915882
$ws.Range("I2").Value = -0.2262087369868314
$ws.Range("G3").Value = 0.0741701026178329
$ws.Range("H3").Value = 8.427783218723121
$ws.Range("G4").Value = -0.03301802634224066
$ws.Range("H4").Value = 26.45178796028735
$ws.Range("G5").Value = -0.03808990438667764
$ws.Range("H5").Value = -40.69934231312891
$ws.Range("G6").Value = -0.09143184112761711
$ws.Range("H6").Value = 13.77705781638189
$ws.Range("G7").Value = -0.08957245007388027
$ws.Range("H7").Value = 1.968583607996397
$ws.Range("G8").Value = -0.3556989430201434
$ws.Range("H8").Value = 3.054087781602297
$ws.Range("G9").Value = -0.3508006414752723
$ws.Range("H9").Value = 10.07203728830305
$ws.Range("G10").Value = 0.01792757505188416
$ws.Range("H10").Value = -11.2100594706433
$ws.Range("G11").Value = 0.02791102977143602
$ws.Range("H11").Value = 22.98600644404992
$ws.Range("G12").Value = 0.2206184719117437
$ws.Range("H12").Value = -0.5091731761476816
$ws.Range("G13").Value = 0.2277134148328923
$ws.Range("H13").Value = 1.113027421822159
$ws.Range("G14").Value = -0.04851302498297466
$ws.Range("H14").Value = -15.21887338943321
$ws.Range("G15").Value = -0.0487118052849469
$ws.Range("H15").Value = -2.129114578782667
$ws.Range("G16").Value = 0.2135892290238797
$ws.Range("H16").Value = 0.4778727053094423
$ws.Range("G17").Value = 0.2207698477432422
$ws.Range("H17").Value = 0.1006054817253016
$ws.Range("G18").Value = 0.07249783660583334
$ws.Range("H18").Value = -0.7154235337864823
$ws.Range("G19").Value = 0.08048572690812149
$ws.Range("H19").Value = 6.83282439896913
$ws.Range("G20").Value = -0.09048711754374515
$ws.Range("H20").Value = -20.67849133658928
$ws.Range("G21").Value = -0.08066702436410382
$ws.Range("H21").Value = 6.817473364614282
$ws.Range("G22").Value = 0.06718170482092298
$ws.Range("H22").Value = -8.598678789851956
$ws.Range("G23").Value = 0.06974205822852109
$ws.Range("H23").Value = 2.06456124413668
$ws.Range("G24").Value = 0.0707439371626922
$ws.Range("H24").Value = 6.202783034296573
$ws.Range("G25").Value = 0.06821457145943334
$ws.Range("H25").Value = 24.52446930104196
$ws.Range("G26").Value = 0.1131917036803619
$ws.Range("H26").Value = -5.159095405303686
$ws.Range("G27").Value = 0.1131917036803619
$ws.Range("H27").Value = -0.5871056892357285
$ws.Range("G28").Value = 0.1374004353846303
$ws.Range("H28").Value = 6.303167120634215
$ws.Range("G29").Value = 0.1424449916137252
$ws.Range("H29").Value = -5.56576406208126
$ws.Range("G30").Value = 0.09010600600179144
$ws.Range("H30").Value = 6.878050420303561
$ws.Range("G31").Value = 0.09010600600179144
$ws.Range("H31").Value = 10.30501011706996
$ws.Range("G32").Value = 0.04892405057673115
$ws.Range("H32").Value = -8.313071724600976
$ws.Range("G33").Value = 0.04982525213316048
$ws.Range("H33").Value = -9.806975446804417
$ws.Range("G34").Value = 0.02238220178775172
$ws.Range("H34").Value = 28.95020453924961
$ws.Range("G35").Value = 0.02619151166756904
$ws.Range("H35").Value = 54.9730487672114
$ws.Range("G36").Value = -0.03096705443833292
$ws.Range("H36").Value = -6.614934068782915
$ws.Range("G37").Value = -0.02653469482251133
$ws.Range("H37").Value = 20.22970296855697
$ws.Range("G38").Value = 0.0806357316957312
$ws.Range("H38").Value = 3.005782900950637
$ws.Range("G39").Value = 0.07944247953606812
$ws.Range("H39").Value = 2.183221750503703
$ws.Range("G40").Value = 0.06936643804804633
$ws.Range("H40").Value = 4.768845492532772
$ws.Range("G41").Value = 0.06919037987643975
$ws.Range("H41").Value = 6.409348555341356
$ws.Range("G42").Value = 0.08319789030811463
$ws.Range("H42").Value = 6.950131515989528
$ws.Range("G43").Value = 0.07505513666577708
$ws.Range("H43").Value = -6.373589166726902
$ws.Range("G44").Value = 0.08859085479851819
$ws.Range("H44").Value = 0.389306752053986
$ws.Range("G45").Value = 0.08822291693893576
$ws.Range("H45").Value = -2.390467887198044
$ws.Range("G46").Value = 0.001916695164717195
$ws.Range("H46").Value = 170.0489785925087
$ws.Range("G47").Value = -0.00019809282802377
$ws.Range("H47").Value = -314.6509217597518
$ws.Range("G48").Value = -0.09979011981959779
$ws.Range("H48").Value = -3.834316661725134
$ws.Range("G49").Value = -0.09839599166905444
$ws.Range("H49").Value = 10.20313819606377
$ws.Range("G50").Value = 0.1635259167591644
$ws.Range("H50").Value = -4.089972543204698
$ws.Range("G51").Value = 0.1788932423670189
$ws.Range("H51").Value = 5.333636029556085
$ws.Range("G52").Value = 0.07296547455520383
$ws.Range("H52").Value = 2.823473439776055
$ws.Range("G53").Value = 0.06417439066526878
$ws.Range("H53").Value = -0.2203451576281396
$ws.Range("G54").Value = -0.1432810377691618
$ws.Range("H54").Value = -12.11098299790892
$ws.Range("G55").Value = -0.1243431967486445
$ws.Range("H55").Value = -6.758779295392643
$ws.Range("G56").Value = 0.194518659283254
$ws.Range("H56").Value = 2.367016309161349
$ws.Range("G57").Value = 0.2088423883435634
$ws.Range("H57").Value = 4.994319122638798
